$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Estimations - update % done ---
$ws.Range("H9").Value = 0.8

# --- Row 10: Scope & out of Scope - set start date, work days, % done ---
# (set days/date before the dependent formula cell is touched so the
#  computed END date / WORK DAYS formulas pick up the final inputs)
$ws.Range("G10").Value = 3
$ws.Range("E10").Value = 43193
$ws.Range("H10").Value = 0.7

# --- Row 11: Size & Effort ---
$ws.Range("G11").Value = 3
$ws.Range("E11").Value = 43193
$ws.Range("H11").Value = 0.8

# --- Row 12: Lifecycle - % done ---
$ws.Range("H12").Value = 1

# --- Row 13: Project Planing - % done (dates/days already set) ---
$ws.Range("H13").Value = 0.8

# --- Row 15: Define the Objective ---
$ws.Range("G15").Value = 7
$ws.Range("E15").Value = 43189
$ws.Range("H15").Value = 0.8

# --- Row 16: Configuration Management Plan ---
$ws.Range("G16").Value = 7
$ws.Range("E16").Value = 43189
$ws.Range("H16").Value = 0.9

# --- Row 17: Test Plan ---
$ws.Range("G17").Value = 7
$ws.Range("E17").Value = 43189
$ws.Range("H17").Value = 0.9

# --- Row 18: Define the Risks ---
$ws.Range("G18").Value = 7
$ws.Range("E18").Value = 43189
$ws.Range("H18").Value = 0.9

# --- Row 22: SRS ---
$ws.Range("G22").Value = 5
$ws.Range("E22").Value = 43189
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 3

# --- Row 23: RTM ---
$ws.Range("G23").Value = 7
$ws.Range("E23").Value = 43189
$ws.Range("I23").Value = 5

# --- Row 24: System Validation Testcases ---
$ws.Range("G24").Value = 4
$ws.Range("E24").Value = 43193

# Force a full recalculation so every dependent formula (END date,
# WORK DAYS, week headers, etc.) is refreshed and cached consistently.
$excel.CalculateFull()

# Restore the cursor to where the author left it.
$ws.Range("Q13").Select()
